# The edit rotates the data held in rows 24-26 while leaving the columns
# that are identical across the three rows untouched (C, S, T, U, V, W, Y,
# Z, AA, AB, AD, AE, AG, AT, AW, AX, AY):
#   new row24 = old row25
#   new row25 = old row26
#   new row26 = old row24

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- snapshot of the "before" values that are needed (rows 24, 25, 26) ----

# Row 24 (before)
$r24 = @{
    A = 111958205; B = 96348; D = "VU"; E = 220787
    F = "Knärot"; G = "Goodyera repens"; H = "(L.) R. Br."; I = "1"
    J = "plantor/tuvor"
    P = "Österåsen, Ång"; Q = 609802.6803741428; R = 7011969.124995505
    AC = "½ m2"
}

# Row 25 (before)
$r25 = @{
    A = 111958182; B = 55611; D = "NT"; E = 102612
    F = "Järpe"; G = "Tetrastes bonasia"; H = "(Linnaeus, 1758)"; I = "3"
    L = "hona"
    P = "Österåsen, Österås, Ång"; Q = 609746.731343443; R = 7011953.229753771
    AC = "1K"
}

# Row 26 (before)
$r26 = @{
    A = 111957798; B = 89686; D = "NT"; E = 658
    F = "Rosenticka"; G = "Rhodofomes roseus"; H = "(Alb. & Schwein.) Kotl. & Pouzar"; I = "6"
    J = "fruktkroppar"
    P = "Österåsen, Österås, Ång"; Q = 609746.731343443; R = 7011953.229753771
}

function Set-TextValue($range, $text) {
    # The "Antal" column stores digit-only values as text. Plain assignment
    # of a numeric-looking string auto-converts to a number, so the cell is
    # forced to Text format first and restored to the default style after.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# ---- write the rotated data back ----

# New row 24 = old row 25
$ws.Range("A24").Value = $r25.A
$ws.Range("B24").Value = $r25.B
$ws.Range("D24").Value = $r25.D
$ws.Range("E24").Value = $r25.E
$ws.Range("F24").Value = $r25.F
$ws.Range("G24").Value = $r25.G
$ws.Range("H24").Value = $r25.H
Set-TextValue $ws.Range("I24") $r25.I
$ws.Range("J24").Value = ""
$ws.Range("L24").Value = $r25.L
$ws.Range("P24").Value = $r25.P
$ws.Range("Q24").Value = $r25.Q
$ws.Range("R24").Value = $r25.R
$ws.Range("AC24").Value = $r25.AC
$ws.Range("AF24").Value = ""

# New row 25 = old row 26
$ws.Range("A25").Value = $r26.A
$ws.Range("B25").Value = $r26.B
$ws.Range("D25").Value = $r26.D
$ws.Range("E25").Value = $r26.E
$ws.Range("F25").Value = $r26.F
$ws.Range("G25").Value = $r26.G
$ws.Range("H25").Value = $r26.H
Set-TextValue $ws.Range("I25") $r26.I
$ws.Range("J25").Value = $r26.J
$ws.Range("L25").Value = ""
$ws.Range("M25").Value = ""
$ws.Range("P25").Value = $r26.P
$ws.Range("Q25").Value = $r26.Q
$ws.Range("R25").Value = $r26.R
$ws.Range("AC25").Value = ""

# New row 26 = old row 24
$ws.Range("A26").Value = $r24.A
$ws.Range("B26").Value = $r24.B
$ws.Range("D26").Value = $r24.D
$ws.Range("E26").Value = $r24.E
$ws.Range("F26").Value = $r24.F
$ws.Range("G26").Value = $r24.G
$ws.Range("H26").Value = $r24.H
Set-TextValue $ws.Range("I26") $r24.I
$ws.Range("J26").Value = $r24.J
$ws.Range("L26").Value = ""
$ws.Range("P26").Value = $r24.P
$ws.Range("Q26").Value = $r24.Q
$ws.Range("R26").Value = $r24.R
$ws.Range("AC26").Value = $r24.AC
